$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H5 with the refined timestamp value
$ws.Range("H5").Value = 46015.30790905093

# Add new row 6 with payment data
$ws.Range("A6").Value = "sample1.PNG"
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = 1042017
$ws.Range("D6").Value = "INR"
$ws.Range("E6").Value = $true
$ws.Range("F6").Value = 0.02
$ws.Range("G6").Value = "READY_FOR_PAYMENT"
$ws.Range("H6").Value = 46017.42774371999
$ws.Range("I6").Value = 0.1215
$ws.Range("J6").Value = 0.786
$ws.Range("K6").Value = ""
$ws.Range("L6").Value = "AUTO_AUDIT_PASSED"
$ws.Range("M6").Value = "OK"

# Match the style of H5 (date/time number format) for H6
$ws.Range("H6").NumberFormat = $ws.Range("H5").NumberFormat
